# Extra wait for loader added in user, ba, customer and withholding page
#
# The underlying xlsx test-data edit regenerates several randomized test
# values in the Customer, BA, User and WithholdingTax sheets. Numeric-looking
# test codes must stay text cells (so leading zeros / shared-string typing is
# preserved) -- force text via NumberFormat "@" then restore a plain style.

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# --- Customer sheet ---
$ws = $wb.Worksheets.Item("Customer")
$ws.Range("A2").Value = "test_cohivg"
Set-TextValue $ws.Range("B2") "391945"
$ws.Range("A4").Value = "test_xahvnh"
Set-TextValue $ws.Range("B4") "467289"

# --- BA sheet ---
$ws = $wb.Worksheets.Item("BA")
Set-TextValue $ws.Range("A2") "366055"
$ws.Range("B2").Value = "test_acrxtu"
Set-TextValue $ws.Range("A4") "366055"
$ws.Range("B4").Value = "test_upzbwk"

# --- User sheet ---
$ws = $wb.Worksheets.Item("User")
$ws.Range("A2").Value = "test_ysamox"
Set-TextValue $ws.Range("B2") "087995"
$ws.Range("A3").Value = "test_afhozv"
Set-TextValue $ws.Range("B3") "654812"
$ws.Range("A5").Value = "test_dcoyki"
Set-TextValue $ws.Range("B5") "137331"

# --- WithholdingTax sheet ---
$ws = $wb.Worksheets.Item("WithholdingTax")
Set-TextValue $ws.Range("B2") "93"
Set-TextValue $ws.Range("C2") "37"
